# Swap the two "Przeniesiono" entries currently on rows 9 and 10 of the
# "Oddziały" sheet: the Wojciechowski (j. angielski) entry moves up to row 9,
# and the Najwer (informatyka) entry moves down to row 10 with an updated
# lesson time and a new remark.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oddziały")

# Capture the current (pre-edit) contents of row 9 and row 10, columns A-G.
$row9 = @()
$row10 = @()
for ($c = 1; $c -le 7; $c++) {
    $row9 += , $ws.Cells.Item(9, $c).Value()
    $row10 += , $ws.Cells.Item(10, $c).Value()
}

# New row 9 becomes the old row 10 (Wojciechowski / j. angielski), unchanged.
for ($c = 1; $c -le 7; $c++) {
    $ws.Cells.Item(9, $c).Value = $row10[$c - 1]
}

# New row 10 becomes the old row 9 (Najwer / informatyka), with the lesson
# time (column B) and remark (column G) updated.
for ($c = 1; $c -le 7; $c++) {
    $ws.Cells.Item(10, $c).Value = $row9[$c - 1]
}
$ws.Range("B10").Value = "18.12.2025, 6, 12:25-13:10, sala: 38"
$ws.Range("G10").Value = "p. Najwer, informatyka za lekcję 3"
